# Apply updated crypto price/volume values (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.150.19"
$ws.Range("E2").Value = "  -4.39%  "
$ws.Range("D3").Value = "'1.651.04"
$ws.Range("E3").Value = "  -3.66%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'215.38"
$ws.Range("E5").Value = "  -4.29%  "
$ws.Range("D6").Value = "'0.5121"
$ws.Range("E6").Value = "  -3.37%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'0.2592"
$ws.Range("E8").Value = "  -2.64%  "
$ws.Range("D9").Value = "'0.06431"
$ws.Range("E9").Value = "  -4.16%  "
$ws.Range("D10").Value = "'19.93"
$ws.Range("E10").Value = "  -4.63%  "
$ws.Range("D11").Value = "'0.07777"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "'1.654.38"
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("D13").Value = "'4.289"
$ws.Range("E13").Value = "  -4.73%  "
$ws.Range("D14").Value = "'1.878.05"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").Value = "'0.5514"
$ws.Range("E15").Value = "  -5.90%  "
$ws.Range("D16").Value = "'0.0₅8002"
$ws.Range("E16").Value = "  -2.76%  "
$ws.Range("D17").Value = "'64.08"
$ws.Range("E17").Value = "  -5.79%  "
$ws.Range("D18").Value = "'26.162.98"
$ws.Range("E18").Value = "  -4.28%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "'210.11"
$ws.Range("E20").Value = "  -5.38%  "
$ws.Range("D21").Value = "'4.396"
$ws.Range("E21").Value = "  -5.62%  "
$ws.Range("E22").Value = "  -4.11%  "
$ws.Range("D23").Value = "'6.055"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").Value = "'143.83"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("D27").Value = "'0.1175"
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("D28").Value = "'6.966"
$ws.Range("E28").Value = "  -3.86%  "
$ws.Range("D29").Value = "'15.80"
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("D30").Value = "'0.05093"
$ws.Range("E30").Value = "  -4.88%  "
$ws.Range("D31").Value = "'1.242"
$ws.Range("E31").Value = "  -4.03%  "
$ws.Range("D32").Value = "'3.352"
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("D33").Value = "'3.216"
$ws.Range("E33").Value = "  -6.27%  "
$ws.Range("D34").Value = "'1.559"
$ws.Range("E34").Value = "  -4.50%  "
$ws.Range("D35").Value = "'2.739"
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("D36").Value = "'0.9238"
$ws.Range("E36").Value = "  -3.63%  "
$ws.Range("D37").Value = "'2.352"
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("D38").Value = "'0.5708"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").Value = "'1.157.14"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").Value = "'0.01586"
$ws.Range("E40").Value = "  -3.47%  "
$ws.Range("D41").Value = "'2.563"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'5.653"
$ws.Range("E43").Value = "  -2.35%  "
$ws.Range("D44").Value = "'0.8241"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("D45").Value = "'100.22"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").Value = "'1.788.58"
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("E47").Value = "  +4.24%  "
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "'55.47"
$ws.Range("E49").Value = "  -3.77%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'7.836"
$ws.Range("E51").Value = "  -3.83%  "
